# Fruta / hortaliza, semanal
#
# A new weekly price record was inserted at row 191 (the sheet's single
# data table is sorted with the most-recent reading on top / near the
# top of this block), pushing all subsequent rows (old 191-249) down by
# one row to (192-250). The newly inserted row carries this week's
# observation; every other row's data is untouched, it just shifts down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 191, shifting rows 191:249 -> 192:250
$ws.Rows.Item(191).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(191, 1).Value  = 3
$ws.Cells.Item(191, 2).Value  = 'Femacal de La Calera'
$ws.Cells.Item(191, 3).Value  = 'Coquimbo'
$ws.Cells.Item(191, 4).Value  = 44627
$ws.Cells.Item(191, 5).Value  = 5
$ws.Cells.Item(191, 6).Value  = 100112001
$ws.Cells.Item(191, 7).Value  = 'Berenjena'
$ws.Cells.Item(191, 8).Value  = 'Sin especificar'
$ws.Cells.Item(191, 9).Value  = 'Primera'
$ws.Cells.Item(191, 10).Value = 50
$ws.Cells.Item(191, 11).Value = 10000
$ws.Cells.Item(191, 12).Value = 10000
$ws.Cells.Item(191, 13).Value = 10000
$ws.Cells.Item(191, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(191, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(191, 16).Value = 167
$ws.Cells.Item(191, 17).Value = 60
$ws.Cells.Item(191, 18).Value = 'Hortaliza'
